$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "C:Temperature" / "temp_mean_cont_scale_clst" row (row 15),
# shifting the rows below it up.
$ws.Rows.Item(15).Delete()

# Update the active selection to B9, matching the saved view state.
$ws.Range("B9").Select()
